$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the marking scheme: correct answer mark 5 -> 4, wrong answer penalty -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Update the totals to reflect the corrected marking scheme
$ws.Range("B12").Value = 100
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "96 / 112"
